$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10: Sending cluster / Ligand symbol / Receptor symbol / Target cluster
# plus the 16 numeric NATMI statistic columns (E..T), recomputed to include the new "ECs" cluster.
$rows = @(
    @("ECs", "Rarres2", "Gpr1", "ECs", 2, 0.6666666666666666, 1.445726333333333, 4.337179, 0.01544456920939864, 0.01544456920939864, 2, 0.6666666666666666, 2.078629666666667, 6.235889, 0.4099944837797483, 0.4099944837797483, 3.005129646347889, 27.046166817131, 0.006332188180207989, 0.00633218818020799),
    @("ECs", "Rarres2", "Gpr1", "FAPs", 2, 0.6666666666666666, 1.445726333333333, 4.337179, 0.01544456920939864, 0.01544456920939864, 3, 1, 2.374602333333333, 7.123807, 0.4683729254179408, 0.4683729254179408, 3.433025124494777, 30.897226120453, 0.007233818062425893, 0.007233818062425893),
    @("ECs", "Rarres2", "Gpr1", "sCs", 2, 0.6666666666666666, 1.445726333333333, 4.337179, 0.01544456920939864, 0.01544456920939864, 3, 1, 0.6166646666666667, 1.849994, 0.1216325908023109, 0.1216325908023109, 0.8915283474362222, 8.023755126926, 0.001878562966764755, 0.001878562966764755),
    @("FAPs", "Rarres2", "Gpr1", "ECs", 3, 1, 75.55280566666666, 226.658417, 0.8071240795570661, 0.8071240795570661, 2, 0.6666666666666666, 2.078629666666667, 6.235889, 0.4099944837797483, 0.4099944837797483, 157.0463032586348, 1413.416729327713, 0.3309164203442038, 0.3309164203442039),
    @("FAPs", "Rarres2", "Gpr1", "FAPs", 3, 1, 75.55280566666666, 226.658417, 0.8071240795570661, 0.8071240795570661, 3, 1, 2.374602333333333, 7.123807, 0.4683729254179408, 0.4683729254179408, 179.4078686259465, 1614.670817633519, 0.3780350663174059, 0.3780350663174059),
    @("FAPs", "Rarres2", "Gpr1", "sCs", 3, 1, 75.55280566666666, 226.658417, 0.8071240795570661, 0.8071240795570661, 3, 1, 0.6166646666666667, 1.849994, 0.1216325908023109, 0.1216325908023109, 46.59074572216645, 419.316711499498, 0.09817259289545646, 0.09817259289545645),
    @("sCs", "Rarres2", "Gpr1", "ECs", 3, 1, 16.60889166666666, 49.82667499999999, 0.1774313512335352, 0.1774313512335352, 2, 0.6666666666666666, 2.078629666666667, 6.235889, 0.4099944837797483, 0.4099944837797483, 34.52373494878611, 310.713614539075, 0.07274587525533646, 0.07274587525533646),
    @("sCs", "Rarres2", "Gpr1", "FAPs", 3, 1, 16.60889166666666, 49.82667499999999, 0.1774313512335352, 0.1774313512335352, 3, 1, 2.374602333333333, 7.123807, 0.4683729254179408, 0.4683729254179408, 39.43951290574721, 354.955616151725, 0.08310404103810903, 0.08310404103810902),
    @("sCs", "Rarres2", "Gpr1", "sCs", 3, 1, 16.60889166666666, 49.82667499999999, 0.1774313512335352, 0.1774313512335352, 3, 1, 0.6166646666666667, 1.849994, 0.1216325908023109, 0.1216325908023109, 10.24211664332778, 92.17904978995, 0.02158143494008969, 0.02158143494008968)
)

$startRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowData[$c]
    }
}

